$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list with refreshed prices / volume(1h) figures,
# and the Cronos / NEARProtocol rows swapped (row 48 <-> row 49).

$ws.Range("D2").Value = "27.107.05"
$ws.Range("E2").Value = "  -0.38%  "
$ws.Range("D3").Value = "1.827.06"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Formula = "'312.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("D7").Formula = "'0.4582"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +7.65%  "
$ws.Range("D8").Formula = "'0.3745"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.12%  "
$ws.Range("D9").Formula = "'0.07333"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.71%  "
$ws.Range("D10").Formula = "'0.8625"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.39%  "
$ws.Range("D11").Formula = "'21.01"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("D12").Value = "1.830.74"
$ws.Range("E12").Value = "  +0.78%  "
$ws.Range("D13").Formula = "'6.712"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.83%  "
$ws.Range("D14").Formula = "'93.03"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.95%  "
$ws.Range("D15").Formula = "'5.365"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.31%  "
$ws.Range("D16").Formula = "'0.07090"
$ws.Range("D16").Style = "Normal"
$ws.Range("E17").Value = "  -0.32%  "
$ws.Range("D18").Formula = "'0.000008862"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("D21").Value = "27.139.18"
$ws.Range("E21").Value = "  -0.36%  "
$ws.Range("D22").Formula = "'5.203"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.22%  "
$ws.Range("E23").Value = "  +1.55%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").Formula = "'151.90"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.75%  "
$ws.Range("D26").Formula = "'2.238"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.30%  "
$ws.Range("E27").Value = "  +1.32%  "
$ws.Range("D28").Formula = "'5.283"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.33%  "
$ws.Range("D29").Formula = "'117.57"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.38%  "
$ws.Range("D30").Formula = "'0.08918"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.41%  "
$ws.Range("D31").Formula = "'0.7661"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.14%  "
$ws.Range("D32").Formula = "'1.199"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.68%  "
$ws.Range("D33").Formula = "'2.974"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.55%  "
$ws.Range("D34").Formula = "'4.480"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.80%  "
$ws.Range("E35").Value = "  -0.26%  "
$ws.Range("E36").Value = "  -0.56%  "
$ws.Range("D37").Formula = "'0.01973"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.73%  "
$ws.Range("D38").Formula = "'0.05299"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.82%  "
$ws.Range("D39").Formula = "'0.5378"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.25%  "
$ws.Range("D40").Formula = "'7.198"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.06%  "
$ws.Range("D41").Formula = "'2.883"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.44%  "
$ws.Range("D42").Formula = "'0.1719"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.70%  "
$ws.Range("D43").Formula = "'0.5207"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +10.93%  "
$ws.Range("D44").Formula = "'8.644"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.60%  "
$ws.Range("D45").Formula = "'10.74"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.84%  "
$ws.Range("D46").Formula = "'1.991"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +11.24%  "
$ws.Range("D47").Formula = "'106.06"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.41%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Formula = "'1.688"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.90%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Formula = "'0.06457"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.14%  "
$ws.Range("D50").Formula = "'1.000"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.28%  "
$ws.Range("D51").Formula = "'0.9264"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.82%  "
